$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.203.86"
$ws.Range("E2").Value = "  -5.14%  "
$ws.Range("D3").Value = "2.332.14"
$ws.Range("E3").Value = "  -6.80%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "508.40"
$ws.Range("E5").Value = "  -4.08%  "
$ws.Range("D6").Value = "125.16"
$ws.Range("E6").Value = "  -7.07%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  -3.78%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.338.48"
$ws.Range("E9").Value = "  -7.64%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.0938"
$ws.Range("E10").Value = "  -5.66%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  -3.22%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "4.71"
$ws.Range("E12").Value = "  -9.99%  "
$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "0.312"
$ws.Range("E13").Value = "  -6.85%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.747.51"
$ws.Range("E14").Value = "  -7.21%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "56.049.78"
$ws.Range("E15").Value = "  -4.89%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "20.96"
$ws.Range("E16").Value = "  -6.79%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000128"
$ws.Range("E17").Value = "  -6.51%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.342.73"
$ws.Range("E18").Value = "  -7.65%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "10.03"
$ws.Range("E19").Value = "  -6.81%  "
$ws.Range("D20").Value = "303.37"
$ws.Range("E20").Value = "  -6.46%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "3.94"
$ws.Range("E21").Value = "  -6.67%  "
$ws.Range("D22").Value = "0.989"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.91"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "62.96"
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.500.03"
$ws.Range("E26").Value = "  -5.71%  "
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").Value = "0.385"
$ws.Range("E27").Value = "  -5.80%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.153"
$ws.Range("E28").Value = "  -5.50%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "7.03"
$ws.Range("E29").Value = "  -6.68%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "171.56"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").Value = "1.64"
$ws.Range("E31").Value = "  -6.22%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0696"
$ws.Range("E32").Value = "  -8.59%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "5.95"
$ws.Range("E33").Value = "  -6.97%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.10"
$ws.Range("E34").Value = "  -11.74%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "17.23"
$ws.Range("E37").Value = "  -5.71%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.15"
$ws.Range("E38").Value = "  -9.36%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "3.63"
$ws.Range("E39").Value = "  -9.64%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "35.54"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("D41").Value = "1.39"
$ws.Range("E41").Value = "  -10.01%  "
$ws.Range("B42").Value = "SuiNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D42").Value = "0.761"
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.25"
$ws.Range("E43").Value = "  -7.27%  "
$ws.Range("D44").Value = "4.64"
$ws.Range("E44").Value = "  -11.27%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.558"
$ws.Range("E45").Value = "  -6.81%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.0895"
$ws.Range("E46").Value = "  -3.16%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "245.72"
$ws.Range("E47").Value = "  -13.10%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "118.04"
$ws.Range("E48").Value = "  -11.77%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0478"
$ws.Range("E49").Value = "  -6.40%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0203"
$ws.Range("E50").Value = "  -7.44%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "16.26"
$ws.Range("E51").Value = "  -8.59%  "
